$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13679
$ws.Range("F2").Value = 13357
$ws.Range("G2").Value = 16085
$ws.Range("H2").Value = 12634
$ws.Range("I2").Value = 19088
$ws.Range("J2").Value = 11254
$ws.Range("K2").Value = 12130
$ws.Range("L2").Value = 11825

$ws.Range("K2").Select()
